$d = $word.ActiveDocument

# Applies the documented edit (todo -> to-do, and "want to" -> "can") while
# preserving the surrounding runs exactly. The underlying engine merges the
# run being text-edited with its immediately preceding sibling whenever they
# share identical run formatting, so a straight Find/Replace on the full
# phrase would swallow the unrelated run that precedes it (e.g. "create ").
# To avoid that we isolate the piece we are about to rename by temporarily
# toggling Bold on the *prefix* span (a pure formatting edit, which this
# engine splits cleanly without merging neighbours), rename the now-isolated
# word, restore the prefix formatting, and finally use the same Bold
# toggle trick to carve the freshly inserted replacement text away from the
# trailing suffix text. The net effect reproduces exactly the 3-way run
# split (prefix / replacement / suffix) that the target document has, with
# the prefix run keeping its original identity and the two new pieces
# appearing as brand new runs (no rsid attributes), matching real Word's
# behaviour when a user retypes part of a sentence.
function Split-And-Rename($fullText, $wordText, $replacement) {
    $r = $d.Content
    $r.Find.ClearFormatting()
    $r.Find.Execute($fullText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $r.Find.Found) {
        Write-Output "NOT FOUND: $fullText"
        return
    }
    $fullStart = $r.Start
    $fullEnd = $r.End

    $wIdx = $fullText.IndexOf($wordText)
    $prefixLen = $wIdx

    $wordStart = $fullStart + $prefixLen
    $wordEnd = $wordStart + $wordText.Length

    # Step 1: bold the prefix to isolate it from the rename edit
    if ($prefixLen -gt 0) {
        $pre = $d.Range($fullStart, $wordStart)
        $pre.Font.Bold = $true
    }

    # Step 2: rename word -> replacement (this will merge with whatever
    # follows it that shares formatting, which is fine - we split it back
    # apart in step 4)
    $sub = $d.Range($wordStart, $wordEnd)
    $sub.Text = $replacement

    # Step 3: restore the prefix formatting
    if ($prefixLen -gt 0) {
        $pre2 = $d.Range($fullStart, $wordStart)
        $pre2.Font.Bold = $false
    }

    # Step 4: split the replacement text away from the trailing suffix text
    $replStart = $wordStart
    $replEnd = $wordStart + $replacement.Length
    $mid = $d.Range($replStart, $replEnd)
    $mid.Font.Bold = $true
    $mid2 = $d.Range($replStart, $replEnd)
    $mid2.Font.Bold = $false
}

Split-And-Rename "and modify a list of todo items" "todo" "to-do"
Split-And-Rename " the list of todo items" "todo" "to-do"
Split-And-Rename "new todo item" "todo" "to-do"
Split-And-Rename " mark a todo item as completed" "todo" "to-do"
Split-And-Rename "to delete a todo item" "todo" "to-do"
Split-And-Rename "So I want to leave feedback or request follow up" "want to" "can"
